$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# MACRO_SCORE (column N) was recomputed; rows 2-5 all share the same score.
$ws.Range("N2:N5").Value = 54.86376272656823
